$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F5").Value = 1866
$ws.Range("F7").Value = 7451
$ws.Range("F15").Value = 3473
$ws.Range("F16").Value = 5825
$ws.Range("F18").Value = 615
$ws.Range("F20").Value = 1191
$ws.Range("F21").Value = 360
$ws.Range("F22").Value = 5865
$ws.Range("F26").Value = 4013
$ws.Range("F28").Value = 672
$ws.Range("F29").Value = 1853
$ws.Range("F33").Value = 163
$ws.Range("F34").Value = 109
$ws.Range("F36").Value = 1109
$ws.Range("F37").Value = 477
$ws.Range("F38").Value = 1820
$ws.Range("F40").Value = 360
$ws.Range("F42").Value = 1034
$ws.Range("F48").Value = 141

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F25").Value = 120

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F5").Value = 438
$ws.Range("F13").Value = 1011
$ws.Range("F14").Value = 1451

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 438
$ws.Range("F8").Value = 7451
$ws.Range("F15").Value = 1011
$ws.Range("F18").Value = 3473
$ws.Range("F21").Value = 615
$ws.Range("F23").Value = 1191
$ws.Range("F24").Value = 360
$ws.Range("F25").Value = 5865
$ws.Range("F28").Value = 672
$ws.Range("F30").Value = 1853
$ws.Range("F34").Value = 163
$ws.Range("F35").Value = 109
$ws.Range("F37").Value = 1109
$ws.Range("F38").Value = 1820
$ws.Range("F40").Value = 360
$ws.Range("F42").Value = 1034
$ws.Range("F43").Value = 120
$ws.Range("F48").Value = 141
